# "ran resolve and classify+summarise steps after changes to mapping file"
# Re-running the pipeline zeroed out the Range Status breakdown (no range
# data produced this run) and the Range Analysis species count, and
# recomputed the High Priority break-up (Range / old-Trend-Different rows
# collapsed away, leaving just Trend New + IUCN with new counts).

$wb = $excel.ActiveWorkbook

# --- "Range Status" sheet: species counts all zeroed, percentage column gone ---
$ws2 = $wb.Worksheets.Item("Range Status")
$ws2.Range("B2").Value = 0
$ws2.Range("C2").ClearContents()
$ws2.Range("B3").Value = 0
$ws2.Range("C3").ClearContents()
$ws2.Range("B4").Value = 0
$ws2.Range("C4").ClearContents()
$ws2.Range("B5").Value = 0
$ws2.Range("C5").ClearContents()
$ws2.Range("B6").Value = 0
$ws2.Range("C6").ClearContents()
$ws2.Range("B7").Value = 0
$ws2.Range("C7").ClearContents()

# --- "Species qualification" sheet: Range Analysis count zeroed ---
$ws4 = $wb.Worksheets.Item("Species qualification")
$ws4.Range("B5").Value = 0

# --- "High Priority break-up" sheet: recomputed breakdown ---
$ws5 = $wb.Worksheets.Item("High Priority break-up")

# Trend New row keeps its label, values recomputed
$ws5.Range("B2").Value = 4
$ws5.Range("C2").Value = 25
$ws5.Range("D2").Value = 4
$ws5.Range("E2").Value = 25

# Old row 3 ("Trend Different") becomes "IUCN" with the former row-5 counts
$ws5.Range("A3").Value = "IUCN"
$ws5.Range("B3").Value = 12
$ws5.Range("C3").Value = 75
$ws5.Range("D3").Value = 12
$ws5.Range("E3").Value = 75

# Old rows 4 (Range) and 5 (IUCN) are dropped entirely
$ws5.Rows("4:5").Delete()
